# Add three new event rows (220-222) to the "Tabelle1" sheet, each with a
# date, event name, location, city and an Instagram hyperlink — mirroring
# the layout of the rows immediately above them (rows 11-219).
#
# xlPasteFormats = -4122 (Excel's PasteSpecial enum) — used to pick up the
# exact cell style (border/fill/number-format/font) already used by the
# other "Link" cells in column E, since Hyperlinks.Add() forces Excel's
# built-in "Hyperlink" style onto the cell and we want it to match its
# neighbours instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$formatDonor = "B219"

$rows = @(
    @{
        Row      = 220
        Date     = 45752
        Event    = "BACKDOOR CULTURE (23Uhr)"
        Location = "Trip Studios (hafen7)"
        City     = "Neuss"
        Link     = "https://www.instagram.com/p/DGgWzW7CR0S/?igsh=OTNwazdxeTA1ZDNu"
    },
    @{
        Row      = 221
        Date     = 45758
        Event    = "TRINITY GOES COLOGNE"
        Location = "Elektroküche"
        City     = "Köln"
        Link     = "https://www.instagram.com/p/DGdv-W_IRqT/?igsh=d3g2cXE0czdodjUx"
    },
    @{
        Row      = 222
        Date     = 45717
        Event    = "DOPAIR"
        Location = "Projekt X"
        City     = "Bochum"
        Link     = "https://www.instagram.com/reel/DGgDReNsaAu/?igsh=aWZzb3ZpNWp6ZzJx"
    }
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Range("A$r").Value = $item.Date

    $ws.Range("B$r").Value = $item.Event
    $ws.Range("B$r").NumberFormat = "@"

    $ws.Range("C$r").Value = $item.Location
    $ws.Range("C$r").NumberFormat = "@"

    $ws.Range("D$r").Value = $item.City
    $ws.Range("D$r").NumberFormat = "@"

    $ws.Range("E$r").Value = $item.Link
    $ws.Hyperlinks.Add($ws.Range("E$r"), $item.Link, "", "", $item.Link)

    # Hyperlinks.Add() just stamped the built-in "Hyperlink" style onto
    # E$r; restore the plain bordered/filled text style used by every
    # other row's Link column.
    $ws.Range($formatDonor).Copy()
    $ws.Range("E$r").PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false

"done"
